$wb = $excel.ActiveWorkbook

# --- Sheet "Test Steps": update Results (column H) for several test step rows ---
$ws1 = $wb.Worksheets.Item("Test Steps")
$ws1.Range("H21").Value = "PASS"
$ws1.Range("H22").Value = "PASS"
$ws1.Range("H25").Value = "PASS"
$ws1.Range("H26").Value = "FAIL"
$ws1.Range("H28").Value = "PASS"
$ws1.Range("H31").Value = "FAIL"

# --- Sheet "Test Cases": flip Runmode (Yes/No) for the Collect/Empty bottles cases,
#     and mark the "Empty collected bottles" case as FAIL ---
$ws2 = $wb.Worksheets.Item("Test Cases")
$ws2.Range("C4").Value = "Yes"
$ws2.Range("C5").Value = "No"
$ws2.Range("D5").Value = "FAIL"

# --- Restore / update selections on each sheet, then make "Test Steps" the active tab ---
$ws3 = $wb.Worksheets.Item("Settings")
$ws3.Range("D6").Select() | Out-Null

$ws2.Range("B4").Select() | Out-Null

$ws1.Range("C36").Select() | Out-Null
$ws1.Activate() | Out-Null
